$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Tinh_trang_don" block (A25:A29) renamed / restructured ---
# Old column list (A25:A30): Tinh_trang_don, ma_hoa_don, ma_san_pham,
#   tinh_trang_don, ma_nhan_vien, thoi_gian_xu_ly
# New column list (A25:A29): lich_su_tinh_trang_don, ma_hoa_don,
#   tinh_trang_don, ma_nhan_vien, thoi_gian_xu_ly   (ma_san_pham removed)
$ws.Range("A25").Value2 = "lich_su_tinh_trang_don"
$ws.Range("A27").Value2 = "tinh_trang_don"
$ws.Range("A28").Value2 = "ma_nhan_vien"
$ws.Range("A29").Value2 = "thoi_gian_xu_ly"
$ws.Range("A30").ClearContents()

# --- "san_pham" table header: The_loai -> the_loai (column J header) ---
$ws.Range("J1").Value2 = "the_loai"

# --- view state: selection moved to G15 ---
$ws.Range("G15").Select()
